$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text cells) ---
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Helper: force a cell to hold literal text (e.g. "0" or "***.*") ---
# using a style-14 (General) reference cell (A15) to restore correct style index after assignment
function Set-TextCell($ref, $text) {
    $cell = $ws.Range($ref)
    $styleRef = $ws.Range("A15")
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $styleRef.Copy()
    $cell.PasteSpecial(-4122)
}

# --- Cells changing between numeric <-> placeholder text type ---
Set-TextCell "D15" "0"
Set-TextCell "E15" "***.*"
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("C17").Value = 1
Set-TextCell "C20" "0"
Set-TextCell "G22" "0"
Set-TextCell "H22" "***.*"
Set-TextCell "D23" "0"
Set-TextCell "E23" "***.*"
Set-TextCell "D26" "0"
Set-TextCell "E26" "***.*"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 180
$ws.Range("M15").Value = 27.272727272727
$ws.Range("N15").Value = -30
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 109
$ws.Range("K16").Value = 2.752293577981
$ws.Range("L16").Value = 30.232558139534
$ws.Range("M16").Value = 31.764705882352
$ws.Range("N16").Value = -80.419580419580
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = 19.277108433734
$ws.Range("L17").Value = 90.384615384615
$ws.Range("M17").Value = 86.792452830188
$ws.Range("N17").Value = 5.319148936170
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("I18").Value = 132
$ws.Range("J18").Value = 79
$ws.Range("K18").Value = 67.088607594936
$ws.Range("L18").Value = -9.589041095890
$ws.Range("M18").Value = 33.333333333333
$ws.Range("N18").Value = -84.982935153583
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -15.789473684210
$ws.Range("I19").Value = 799
$ws.Range("J19").Value = 718
$ws.Range("K19").Value = 11.281337047353
$ws.Range("L19").Value = 69.279661016949
$ws.Range("M19").Value = 21.613394216133
$ws.Range("N19").Value = -54.185779816513
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = -18.604651162790
$ws.Range("N20").Value = -93.794326241134
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -27.272727272727
$ws.Range("F21").Value = 74
$ws.Range("H21").Value = -17.777777777777
$ws.Range("I21").Value = 1226
$ws.Range("J21").Value = 1082
$ws.Range("K21").Value = 13.308687615526
$ws.Range("L21").Value = 52.298136645962
$ws.Range("M21").Value = 30.564430244941
$ws.Range("N21").Value = -72.374943668319
$ws.Range("F22").Value = 2
$ws.Range("M22").Value = -21.875
$ws.Range("M23").Value = -4
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 42.857142857142
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = -11.111111111111
$ws.Range("I24").Value = 1194
$ws.Range("J24").Value = 1217
$ws.Range("K24").Value = -1.889893179950
$ws.Range("L24").Value = -8.505747126436
$ws.Range("M24").Value = 17.984189723320
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 38.461538461538
$ws.Range("I25").Value = 220
$ws.Range("J25").Value = 187
$ws.Range("K25").Value = 17.647058823529
$ws.Range("L25").Value = 66.666666666666
$ws.Range("M25").Value = -19.413919413919
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 19
$ws.Range("K26").Value = 35.714285714285
$ws.Range("L26").Value = 137.5
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 53
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = 1.923076923076
$ws.Range("L27").Value = 47.222222222222

Write-Host "edit complete"
